$wb = $excel.ActiveWorkbook

# --- "plans" sheet: data edits -------------------------------------------------
$plans = $wb.Worksheets.Item("plans")

# Column A (id) rows 2-6: was text "planN", now plain integers 1-5
$plans.Cells.Item(2,1).Value = 1
$plans.Cells.Item(3,1).Value = 2
$plans.Cells.Item(4,1).Value = 3
$plans.Cells.Item(5,1).Value = 4
$plans.Cells.Item(6,1).Value = 5

# Row 5 gains expand_by / patch_dim0 / patch_dim1 values
$plans.Cells.Item(5,3).Value = 0
$plans.Cells.Item(5,11).Value = "128"
$plans.Cells.Item(5,12).Value = "96"

# New row 7 (id=6), a copy-like row similar to row 3 (lbd / nodesthick,nodes)
$plans.Cells.Item(7,1).Value = 6
$plans.Cells.Item(7,2).Value = "nodesthick,nodes"
$plans.Cells.Item(7,3).Value = 0
$plans.Cells.Item(7,4).Value = "0"
$plans.Cells.Item(7,5).Value = "0"
$plans.Cells.Item(7,10).Value = "lbd"
$plans.Cells.Item(7,11).Value = "128"
$plans.Cells.Item(7,12).Value = "96"
$plans.Cells.Item(7,13).Value = "0.25"
$plans.Cells.Item(7,16).Value = "TSL.label_localiser,TSL.label_localiser"
$plans.Cells.Item(7,18).Value = 2
$plans.Cells.Item(7,20).Value = "0.8,.8,1.5"
$plans.Cells.Item(7,21).Value = "manual_value"

# V7 flips from FALSE() to TRUE()
$plans.Cells.Item(7,22).Formula = "=TRUE()"

# Best-effort: widen column V (and the sheet's default width) ~1.5x
$plans.Range("V:V").ColumnWidth = 15.11

# --- Selections across the workbook (mirrors clicking around before the edit) --
$wb.Worksheets.Item("dataset_params").Range("E13").Select()
$wb.Worksheets.Item("transform_factors").Range("A1").Select()
$wb.Worksheets.Item("affine3d").Range("A1").Select()
$wb.Worksheets.Item("loss_params").Range("A1").Select()
$wb.Worksheets.Item("plan1").Range("A1").Select()
$wb.Worksheets.Item("plan2").Range("A1").Select()
$wb.Worksheets.Item("plan3").Range("A1").Select()
$wb.Worksheets.Item("plan4").Range("A1").Select()

# "plans" becomes the active/selected tab with B8 selected
$plans.Range("B8").Select()
